$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: clear its text (this also absorbs/removes the old blank
# spacer paragraph that followed it), then mark it as a no-bullet, zero
# indent paragraph (marL="0" indent="0" + buNone).
$para1 = $tr.Paragraphs(1)
$chars = $para1.Characters(1, $para1.Length)
$chars.Delete()

$bullet1 = $para1.ParagraphFormat.Bullet
$bullet1.Visible = 0

$levels = $tf.Ruler.Levels
$lvl1 = $levels.Item(1)
$lvl1.LeftMargin = 0
$lvl1.FirstMargin = 0

# --- Paragraph 2 (previously "Charter schools perform on average worse...
# Hamilton County, though this is not a statistically significant result")
# gets its text replaced. Blank it first so PowerPoint doesn't keep a
# shared-prefix run split between old/new text.
$para2 = $tr.Paragraphs(2)
$para2.Text = " "
$para2.Text = "Charter schools perform on average worse than public schools in Hamilton County, Cuyahoga County, and Franklin County"

# --- Reset the body autofit back to a plain normAutofit (drop the
# fontScale/lnSpcReduction overrides).
$tf.AutoSize = 2
